$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '57.000.66'
$ws.Range('E2').Value = '  -2.34%  '
$ws.Range('D3').Value = '3.062.87'
$ws.Range('E3').Value = '  -2.75%  '
$ws.Range('E4').Value = '  -0.01%  '
$ws.Range('D5').Value = '''520.71'
$ws.Range('E5').Value = '  -2.15%  '
$ws.Range('D6').Value = '''135.02'
$ws.Range('E6').Value = '  -5.48%  '
$ws.Range('E7').Value = '  +0.03%  '
$ws.Range('D8').Value = '3.060.56'
$ws.Range('E8').Value = '  -2.67%  '
$ws.Range('D9').Value = '''0.464'
$ws.Range('E9').Value = '  +3.89%  '
$ws.Range('D10').Value = '''7.27'
$ws.Range('E10').Value = '  +1.50%  '
$ws.Range('E11').Value = '  -3.48%  '
$ws.Range('D12').Value = '''0.401'
$ws.Range('E12').Value = '  +1.74%  '
$ws.Range('E13').Value = '  +1.08%  '
$ws.Range('D14').Value = '3.585.17'
$ws.Range('E14').Value = '  -2.63%  '
$ws.Range('D15').Value = '''25.05'
$ws.Range('E15').Value = '  -2.63%  '
$ws.Range('D16').Value = '''0.0000159'
$ws.Range('E16').Value = '  -4.11%  '
$ws.Range('D17').Value = '57.027.56'
$ws.Range('E17').Value = '  -2.37%  '
$ws.Range('D18').Value = '3.055.39'
$ws.Range('E18').Value = '  -2.49%  '
$ws.Range('D19').Value = '''5.84'
$ws.Range('E19').Value = '  -4.70%  '
$ws.Range('E20').Value = '  -3.79%  '
$ws.Range('D21').Value = '''7.76'
$ws.Range('E21').Value = '  -3.15%  '
$ws.Range('D22').Value = '''347.61'
$ws.Range('E22').Value = '  +0.89%  '
$ws.Range('E23').Value = '  -0.20%  '
$ws.Range('D24').Value = '''68.91'
$ws.Range('E24').Value = '  +1.74%  '
$ws.Range('D25').Value = '''0.496'
$ws.Range('E25').Value = '  -3.73%  '
$ws.Range('E26').Value = '  +0.24%  '
$ws.Range('D27').Value = '''0.163'
$ws.Range('E27').Value = '  -3.92%  '
$ws.Range('D28').Value = '0.0₃0855'
$ws.Range('E28').Value = '  -8.99%  '
$ws.Range('D29').Value = '''0.999'
$ws.Range('E29').Value = '  +0.01%  '
$ws.Range('E30').Value = '  -3.78%  '
$ws.Range('D31').Value = '''1.85'
$ws.Range('E31').Value = '  -2.19%  '
$ws.Range('B32').Value = 'EthereumClassic'
$ws.Range('C32').Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
$ws.Range('D32').Value = '''20.83'
$ws.Range('E32').Value = '  -1.72%  '
$ws.Range('B33').Value = 'RenderToken'
$ws.Range('C33').Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range('D33').Value = '''5.76'
$ws.Range('E33').Value = '  -10.69%  '
$ws.Range('D34').Value = '''158.27'
$ws.Range('E34').Value = '  -0.34%  '
$ws.Range('E35').Value = '  -0.83%  '
$ws.Range('D36').Value = '''1.12'
$ws.Range('E36').Value = '  -6.52%  '
$ws.Range('E37').Value = '  -4.49%  '
$ws.Range('D38').Value = '''25.16'
$ws.Range('E38').Value = '  -4.75%  '
$ws.Range('E39').Value = '  -2.74%  '
$ws.Range('E40').Value = '  -2.57%  '
$ws.Range('D41').Value = '''1.56'
$ws.Range('E41').Value = '  -6.68%  '
$ws.Range('D42').Value = '''4.01'
$ws.Range('E42').Value = '  +0.00%  '
$ws.Range('D43').Value = '''0.690'
$ws.Range('E43').Value = '  -2.41%  '
$ws.Range('D44').Value = '2.395.29'
$ws.Range('E44').Value = '  +5.40%  '
$ws.Range('B45').Value = 'OKB'
$ws.Range('C45').Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range('D45').Value = '''36.37'
$ws.Range('E45').Value = '  -1.00%  '
$ws.Range('B46').Value = 'FirstDigitalUSD'
$ws.Range('C46').Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range('D46').Value = '''0.999'
$ws.Range('E46').Value = '  -0.09%  '
$ws.Range('D47').Value = '3.099.39'
$ws.Range('E47').Value = '  -2.73%  '
$ws.Range('E48').Value = '  -1.81%  '
$ws.Range('D49').Value = '''5.95'
$ws.Range('E49').Value = '  -2.88%  '
$ws.Range('D50').Value = '''0.926'
$ws.Range('E50').Value = '  -8.15%  '
$ws.Range('E51').Value = '  -6.88%  '
